$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.324.19"
$ws.Range("D3").Value = "1.803.76"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'227.74"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "'0.576"
$ws.Range("E6").Value = "  +3.68%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +9.60%  "
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("D10").Value = "'0.0693"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").Value = "'0.0968"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").Value = "2.063.98"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("E13").Value = "  +4.66%  "
$ws.Range("D14").Value = "1.820.35"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "'4.49"
$ws.Range("E16").Value = "  +4.45%  "
$ws.Range("D17").Value = "34.315.72"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'69.07"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "'245.22"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'11.49"
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'4.18"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'172.64"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").Value = "'7.93"
$ws.Range("E26").Value = "  +8.27%  "
$ws.Range("D27").Value = "'16.85"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").Value = "'0.118"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'0.0531"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").Value = "'3.83"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "1.393.06"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("E37").Value = "  -5.81%  "
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("E40").Value = "  +10.73%  "
$ws.Range("D41").Value = "'0.961"
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").Value = "'81.65"
$ws.Range("E43").Value = "  -3.03%  "
$ws.Range("D44").Value = "'2.41"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("D46").Value = "'6.03"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "'0.0501"
$ws.Range("E47").Value = "  -4.95%  "
$ws.Range("D48").Value = "1.964.48"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "'104.85"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  -0.07%  "
